# Apply the "cryptos list" refresh captured by the commit diff.
# Every Price (D) / Volume(1h) (E) cell is stored as literal text in the
# source workbook (t="inlineStr"), and a few rows had their Coin/Link/Price/
# Volume swapped with a neighboring row. We reproduce both kinds of edits
# by writing plain cell values; numeric-looking Price strings are forced to
# stay text (NumberFormat "@" + quote-prefix) and then restored to the
# workbook default style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.270.68"
$ws.Range("E2").Value = "  -0.47%  "

# Row 3
$ws.Range("D3").Value = "3.310.39"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "1.00"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "189.89"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E5").Value = "  +3.68%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "562.19"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E6").Value = "  +0.33%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "0.589"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E8").Value = "  -1.77%  "

# Row 9
$ws.Range("D9").Value = "3.303.61"
$ws.Range("E9").Value = "  -1.71%  "

# Row 10
$ws.Range("E10").Value = "  -1.08%  "

# Row 11
$ws.Range("E11").Value = "  -1.33%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "48.01"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E12").Value = "  -0.21%  "

# Row 13
$ws.Range("E13").Value = "  +1.13%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "8.71"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E14").Value = "  -0.41%  "

# Row 15
$ws.Range("D15").Value = "3.840.83"
$ws.Range("E15").Value = "  -1.89%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "611.54"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E16").Value = "  +0.90%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.289.93"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D18")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "18.09"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E18").Value = "  -1.16%  "

# Row 19
$ws.Range("E19").Value = "  -0.14%  "

# Row 20
$ws.Range("D20").Value = "3.309.65"
$ws.Range("E20").Value = "  -1.91%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "11.13"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E21").Value = "  -4.11%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "0.913"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E22").Value = "  -0.44%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "18.49"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E23").Value = "  +8.79%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "5.12"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E24").Value = "  -0.99%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "101.60"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E25").Value = "  +3.05%  "

# Row 26
$ws.Range("E26").Value = "  -1.59%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "2.76"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E28").Value = "  +1.16%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "9.85"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E29").Value = "  +3.96%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "8.65"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "30.37"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E31").Value = "  -1.33%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "6.78"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E32").Value = "  +7.07%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "4.06"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E33").Value = "  +6.48%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "575.34"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E34").Value = "  +3.55%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "11.13"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E35").Value = "  -0.60%  "

# Row 36
$ws.Range("E36").Value = "  -0.41%  "

# Row 37
$ws.Range("D37").Value = "3.732.36"
$ws.Range("E37").Value = "  -2.01%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D38")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "57.23"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E38").Value = "  -1.65%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D39")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "1.00"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0732"
$ws.Range("E40").Value = "  +1.00%  "

# Row 41
$ws.Range("B41").Value = "CoreDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$c = $ws.Range("D41")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "3.49"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E41").Value = "  +3.37%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D42")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "3.33"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E42").Value = "  -2.38%  "

# Row 43
$ws.Range("E43").Value = "  +5.26%  "

# Row 44
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "2.73"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E45").Value = "  +1.11%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "0.341"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E46").Value = "  -2.45%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "0.0427"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "3.27"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E48").Value = "  +1.45%  "

# Row 49
$ws.Range("E49").Value = "  -1.02%  "

# Row 50
$ws.Range("E50").Value = "  -2.85%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"   # force text so "48.01" etc. is not read back as a number
$c.Value = "1.00"
$c.Style = "Normal"     # drop the quote-prefix style again, keep value as text
$ws.Range("E51").Value = "  +0.11%  "
